# Weekly data refresh: a new week's record is inserted as row 22
# (pushing the existing rows 22-37 down to 23-38, each keeping its
# original values), and the new row 22 is populated with the latest
# week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 22 (row 21 and above are untouched,
# everything from the old row 22 onward shifts down by one).
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with this week's record.
$ws.Range("A22").Value = 9
$ws.Range("B22").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C22").Value = "Metropolitana"
$ws.Range("D22").Value = 44957
$ws.Range("E22").Value = 13
$ws.Range("F22").Value = 100112010
$ws.Range("G22").Value = "Achicoria"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 70
$ws.Range("K22").Value = 7000
$ws.Range("L22").Value = 7000
$ws.Range("M22").Value = 7000
$ws.Range("N22").Value = "`$/caja 16 unidades"
$ws.Range("O22").Value = "Provincia de Quillota"
$ws.Range("P22").Value = 438
$ws.Range("Q22").Value = 16
$ws.Range("R22").Value = "Hortaliza"
